$d = $word.ActiveDocument

# --- Change 1 -----------------------------------------------------------
# The bullet about JWT/HttpOnly cookies was split across three runs so
# that "HttpOnly" could be wrapped in a <w:proofErr> spellStart/spellEnd
# pair. Collapse the whole sentence back down to a single plain run.
$target1 = "Verify authentication/authorization protects admin routes using JWT stored in HttpOnly cookies."

$rng1 = $d.Content
$found1 = $rng1.Find.Execute("JWT stored in*cookies.", $false, $false, $true, $false, $false, `
                              $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not locate the JWT/HttpOnly sentence to update."
}
$rng1.Expand(4) | Out-Null        # wdParagraph -> whole bullet paragraph
$rng1.MoveEnd(1, -1) | Out-Null   # wdCharacter -> drop the trailing paragraph mark

# Re-assigning identical text is treated as a no-op (it would leave the
# original multi-run / proofErr structure untouched), so stamp a
# placeholder first and then overwrite it with the real sentence; the
# range tracks the edit automatically.
$rng1.Text = "TEMP_PLACEHOLDER__"
$rng1.Text = $target1

# --- Change 2 -------------------------------------------------------------
# Remove the whole "Regression Smoke Tests before demos/interview" bullet
# (paragraph text + its paragraph mark) that used to follow the
# "Security Smoke Testing ..." bullet.
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("Regression Smoke Tests before demos/interview", $true, $true, `
                              $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not locate the Regression Smoke Tests paragraph to remove."
}
$rng2.Expand(4) | Out-Null        # wdParagraph -> whole bullet paragraph + mark
$rng2.Delete()
